# Auto-generated cell-value update script
# Applies the numeric corrections captured in the commit diff for Masamune_Profits
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5153.8623
$ws.Range("I40").Value = 7691.3125
$ws.Range("J40").Value = 2030.8462
$ws.Range("K40").Value = 7691.3125
$ws.Range("L40").Value = 2030.8462
$ws.Range("M40").Value = -7516.3125
$ws.Range("N40").Value = -2380.8462
$ws.Range("H86").Value = 3848646
$ws.Range("I86").Value = 6252460.5
$ws.Range("J86").Value = 2543
$ws.Range("K86").Value = 6252460.5
$ws.Range("L86").Value = 2543
$ws.Range("M86").Value = -6251337.5
$ws.Range("N86").Value = -4789
$ws.Range("H89").Value = 3848646
$ws.Range("I89").Value = 6252460.5
$ws.Range("J89").Value = 2543
$ws.Range("K89").Value = 31262302.5
$ws.Range("L89").Value = 12715
$ws.Range("M89").Value = -31256686.5
$ws.Range("N89").Value = -23947
$ws.Range("H132").Value = 24441.205
$ws.Range("I132").Value = 3482
$ws.Range("J132").Value = 85222.89999999999
$ws.Range("K132").Value = 10446
$ws.Range("L132").Value = 255668.7
$ws.Range("M132").Value = -7916
$ws.Range("N132").Value = -260728.7
$ws.Range("H135").Value = 26317054
$ws.Range("I135").Value = 1372.2941
$ws.Range("K135").Value = 12350.6469
$ws.Range("M135").Value = -9815.6469
$ws.Range("H137").Value = 3350163.5
$ws.Range("I137").Value = 6994014.5
$ws.Range("J137").Value = 9966.833000000001
$ws.Range("K137").Value = 20982043.5
$ws.Range("L137").Value = 29900.499
$ws.Range("M137").Value = -20979493.5
$ws.Range("N137").Value = -35000.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11603.13
$ws.Range("I32").Value = 10215.408
$ws.Range("K32").Value = 10215.408
$ws.Range("M32").Value = -9928.407999999999
$ws.Range("H61").Value = 3111.9546
$ws.Range("I61").Value = 2499.75
$ws.Range("J61").Value = 3461.7856
$ws.Range("K61").Value = 2499.75
$ws.Range("L61").Value = 3461.7856
$ws.Range("M61").Value = -2287.75
$ws.Range("N61").Value = -3885.7856
$ws.Range("H136").Value = 3111.9546
$ws.Range("I136").Value = 2499.75
$ws.Range("J136").Value = 3461.7856
$ws.Range("K136").Value = 7499.25
$ws.Range("L136").Value = 10385.3568
$ws.Range("M136").Value = -4949.25
$ws.Range("N136").Value = -15485.3568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 51499.5
$ws.Range("I112").Value = 50000
$ws.Range("J112").Value = 52999
$ws.Range("K112").Value = 50000
$ws.Range("L112").Value = 52999
$ws.Range("M112").Value = -48523
$ws.Range("N112").Value = -55953
$ws.Range("H130").Value = 54969.75
$ws.Range("J130").Value = 54969.75
$ws.Range("L130").Value = 54969.75
$ws.Range("N130").Value = -65009.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5957516
$ws.Range("I31").Value = 1768.4242
$ws.Range("J31").Value = 14502719
$ws.Range("K31").Value = 1768.4242
$ws.Range("L31").Value = 14502719
$ws.Range("M31").Value = -1473.4242
$ws.Range("N31").Value = -14503309
$ws.Range("H34").Value = 5957516
$ws.Range("I34").Value = 1768.4242
$ws.Range("J34").Value = 14502719
$ws.Range("K34").Value = 1768.4242
$ws.Range("L34").Value = 14502719
$ws.Range("M34").Value = -1566.4242
$ws.Range("N34").Value = -14503123
$ws.Range("H62").Value = 3114.7222
$ws.Range("I62").Value = 3298.5
$ws.Range("J62").Value = 2885
$ws.Range("K62").Value = 3298.5
$ws.Range("L62").Value = 2885
$ws.Range("M62").Value = -2674.5
$ws.Range("N62").Value = -4133
$ws.Range("H65").Value = 3114.7222
$ws.Range("I65").Value = 3298.5
$ws.Range("J65").Value = 2885
$ws.Range("K65").Value = 16492.5
$ws.Range("L65").Value = 14425
$ws.Range("M65").Value = -13372.5
$ws.Range("N65").Value = -20665
$ws.Range("H122").Value = 158429.72
$ws.Range("I122").Value = 276126.5
$ws.Range("J122").Value = 1500.6666
$ws.Range("K122").Value = 828379.5
$ws.Range("L122").Value = 4501.9998
$ws.Range("M122").Value = -825929.5
$ws.Range("N122").Value = -9401.9998
$ws.Range("H132").Value = 825797.5600000001
$ws.Range("I132").Value = 1424
$ws.Range("J132").Value = 2337149.2
$ws.Range("K132").Value = 4272
$ws.Range("L132").Value = 7011447.600000001
$ws.Range("M132").Value = -1742
$ws.Range("N132").Value = -7016507.600000001
$ws.Range("H134").Value = 536134.0600000001
$ws.Range("I134").Value = 672193.5
$ws.Range("J134").Value = 178978
$ws.Range("K134").Value = 2016580.5
$ws.Range("L134").Value = 536934
$ws.Range("M134").Value = -2014045.5
$ws.Range("N134").Value = -542004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3219.4614
$ws.Range("J113").Value = 969.8333
$ws.Range("L113").Value = 2909.4999
$ws.Range("N113").Value = -7249.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1127.125
$ws.Range("I102").Value = 980.6
$ws.Range("J102").Value = 1371.3334
$ws.Range("K102").Value = 980.6
$ws.Range("L102").Value = 1371.3334
$ws.Range("M102").Value = 641.4
$ws.Range("N102").Value = -4615.3334
$ws.Range("H107").Value = 4055.6
$ws.Range("I107").Value = 847.5
$ws.Range("J107").Value = 4857.625
$ws.Range("K107").Value = 847.5
$ws.Range("L107").Value = 4857.625
$ws.Range("M107").Value = 1072.5
$ws.Range("N107").Value = -8697.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("H122").Value = 102260.6
$ws.Range("I122").Value = 168667.5
$ws.Range("J122").Value = 2650.25
$ws.Range("K122").Value = 506002.5
$ws.Range("L122").Value = 7950.75
$ws.Range("M122").Value = -503552.5
$ws.Range("N122").Value = -12850.75
$ws.Range("H132").Value = 4213
$ws.Range("I132").Value = 2962
$ws.Range("K132").Value = 8886
$ws.Range("M132").Value = -6356

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3196.5715
$ws.Range("I113").Value = 4344
$ws.Range("J113").Value = 1666.6666
$ws.Range("K113").Value = 13032
$ws.Range("L113").Value = 4999.9998
$ws.Range("M113").Value = -10862
$ws.Range("N113").Value = -9339.9998
$ws.Range("H136").Value = 361462.03
$ws.Range("I136").Value = 534248.3
$ws.Range("K136").Value = 1602744.9
$ws.Range("M136").Value = -1600194.9

# LTW row 109 loses its N109 cell entirely (it previously held -35198;
# the updated source data has no profit figure for this leve/item combo)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N109").ClearContents()

